$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Shift rows 17-19 down to 18-20 (formats first, then literal
#    values), working bottom-up so a source row is never clobbered
#    before it has been used.
# ------------------------------------------------------------------

# old row19 (blank, bordered only) -> row20
$ws.Range("A19:P19").Copy() | Out-Null
$ws.Range("A20:P20").PasteSpecial(-4122) | Out-Null

# old row18 ("SMB" mapping row) -> row19
$ws.Range("A18:P18").Copy() | Out-Null
$ws.Range("A19:P19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = "ZŁOŻENIE"
$ws.Range("B19").Value = "SMB"
$ws.Range("C19").Value = "SBM"
$ws.Range("D19").Value = "Beistellung SBM"
$ws.Range("E19").Value = "Centrum kompletacji"
$ws.Range("F19").Value = "Kompletacja"

# old row17 ("SBM" mapping row) -> row18
$ws.Range("A17:P17").Copy() | Out-Null
$ws.Range("A18:P18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "ZŁOŻENIE"
$ws.Range("B18").Value = "SBM"
$ws.Range("C18").Value = "SBM"
$ws.Range("D18").Value = "Beistellung SBM"
$ws.Range("E18").Value = "Centrum kompletacji"
$ws.Range("F18").Value = "Kompletacja"

# ------------------------------------------------------------------
# 2) Rebuild row 17 as a brand-new "BC" mapping row, formatted like
#    row 16 directly above it.
# ------------------------------------------------------------------
$ws.Range("A16:P16").Copy() | Out-Null
$ws.Range("A17:P17").PasteSpecial(-4122) | Out-Null

$ws.Range("A17").Value = "ZŁOŻENIE"
$ws.Range("B17").Value = "BC"
$ws.Range("C17").Value = "Czarny_Montaż"
$ws.Range("D17").Value = "Spawanie"
$ws.Range("E17").Value = "Ocynkownia"
$ws.Range("F17").Value = "Ocynk"
$ws.Range("G17").Value = "Centrum kompletacji"
$ws.Range("H17").Value = "Kompletacja"
$ws.Range("I17").Value = "Montaż"
$ws.Range("J17").Value = "Biały_Montaż"
$ws.Range("K17").Value = "Odbiór Końcowy"
$ws.Range("L17").Value = "Odbiory"

# ------------------------------------------------------------------
# 3) Fill in the new trailing row 20.
# ------------------------------------------------------------------
$ws.Range("C16:D16").Copy() | Out-Null
$ws.Range("C20:D20").PasteSpecial(-4122) | Out-Null

$ws.Range("A20").Value = "ZŁOŻENIE"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = "Centrum kompletacji"
$ws.Range("D20").Value = "Kompletacja"

# ------------------------------------------------------------------
# 4) Selection / view bookkeeping to match the saved workbook state.
# ------------------------------------------------------------------
$ws.Range("C20:D20").Select() | Out-Null
